$d = $word.ActiveDocument
$dash = [char]8211

# ---------------------------------------------------------------------------
# 1. First paragraph: "This is a Microsoft word document." ->
#    "This is a Microsoft word document.  " (black) +
#    "(This is a change <en-dash> Ve" / "rsion for main branch" / ")" (red)
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("This is a Microsoft word document.", $false, $false, $false, $false, $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

# Position right after the two trailing spaces but still inside paragraph 1
# (position 36 sits between the 2nd trailing space and the paragraph mark).
$ip1 = $d.Range(36, 36)
$ip1.InsertAfter("(This is a change " + $dash + " Ve")
$ip1.Font.Color = 255

$ip2 = $d.Range($ip1.End, $ip1.End)
$ip2.InsertAfter("rsion for main branch")
$ip2.Font.Color = 255

$ip3 = $d.Range($ip2.End, $ip2.End)
$ip3.InsertAfter(")")
$ip3.Font.Color = 255

# ---------------------------------------------------------------------------
# 2. Remove the trailing "ank God almighty, we are free at last." paragraph
#    (keeps "...Shall be lifted-nevermore!" as the final paragraph).
# ---------------------------------------------------------------------------
$lastIdx = $d.Paragraphs.Count
$d.Paragraphs($lastIdx).Range.Delete()

# ---------------------------------------------------------------------------
# 3. Drop the unused styles that Word stripped out on save (none of them are
#    referenced by any paragraph/run in the document). Deleting from the end
#    of the styles collection backwards keeps the runtime's internal indices
#    consistent.
# ---------------------------------------------------------------------------
$d.Styles("podcast-tools__subscribe-links").Delete()
$d.Styles("generic-title").Delete()
$d.Styles("subscribe-more-info").Delete()
$d.Styles("subscribe").Delete()
$d.Styles("audio-tool").Delete()
$d.Styles("Heading 4 Char").Delete()
$d.Styles("Heading 2 Char").Delete()
$d.Styles("Hyperlink").Delete()
$d.Styles("apple-converted-space").Delete()
$d.Styles("Heading 4").Delete()
$d.Styles("Heading 2").Delete()

Write-Output "done"
